# Insert a new claim row above the current row 40 ("-416" / Paraguay 3765 / ...),
# pushing the existing rows 40-49 down to 41-50, then populate the newly
# inserted row 40 with the new claim's data.
#
# Columns A, B, D, E, I hold numeric/date-*looking* text (e.g. "-419",
# "5/26/2025", "13", "806926532", "1") that must stay plain text (matching
# every other row in the sheet), so each of those cells is pre-formatted as
# Text before the value is written to stop Excel auto-converting it to a
# real number/date. M and N are genuine numbers and are written as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 40:49 down to 41:50, creating a blank row 40.
$ws.Rows(40).Insert()

$newRow = 40

$textCols = @(1, 2, 4, 5, 9)   # A, B, D, E, I -> numeric/date-looking text
foreach ($col in $textCols) {
    $ws.Cells.Item($newRow, $col).NumberFormat = "@"
}

$ws.Cells.Item($newRow, 1).Value  = "-419"
$ws.Cells.Item($newRow, 2).Value  = "5/26/2025"
$ws.Cells.Item($newRow, 3).Value  = "Juana Azurduy 1520"
$ws.Cells.Item($newRow, 4).Value  = "13"
$ws.Cells.Item($newRow, 5).Value  = "806926532"
$ws.Cells.Item($newRow, 6).Value  = "GESTION TELECENTRO"
$ws.Cells.Item($newRow, 7).Value  = "Pendiente"
$ws.Cells.Item($newRow, 8).Value  = "Aplomar columna 168 con rienda a pique "
$ws.Cells.Item($newRow, 9).Value  = "1"
$ws.Cells.Item($newRow, 10).Value = "Cambio"
$ws.Cells.Item($newRow, 11).Value = "Fuente TLC"
$ws.Cells.Item($newRow, 12).Value = "Pasante"
$ws.Cells.Item($newRow, 13).Value = -58.458478
$ws.Cells.Item($newRow, 14).Value = -34.546285
